$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K") values per regen of save_data (K instead of Strike#)
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 2
$ws.Range("G7").Value = 1
